$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: new cells D1:L1 need same style as C1 (bold/border style index 1)
$ws.Range("C1").Copy($ws.Range("D1:L1"))
$ws.Range("D1").Value = 2
$ws.Range("E1").Value = 3
$ws.Range("F1").Value = 4
$ws.Range("G1").Value = 5
$ws.Range("H1").Value = 6
$ws.Range("I1").Value = 7
$ws.Range("J1").Value = 8
$ws.Range("K1").Value = 9
$ws.Range("L1").Value = 10

# Row 2: plain numeric cells, no special style (matches C2 which has no s attr)
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 4
$ws.Range("G2").Value = 5
$ws.Range("H2").Value = 6
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
$ws.Range("K2").Value = 12
$ws.Range("L2").Value = 12

# Row 3: C3 value changes, D3:L3 new
$ws.Range("C3").Value = 42.43
$ws.Range("D3").Value = 55.9
$ws.Range("E3").Value = 49.24
$ws.Range("F3").Value = 43.01
$ws.Range("G3").Value = 33.54
$ws.Range("H3").Value = 40.31
$ws.Range("I3").Value = 47.17
$ws.Range("J3").Value = 31.62
$ws.Range("K3").Value = 40.31
$ws.Range("L3").Value = 44.72

# Row 4: C4 value changes, D4:L4 new
$ws.Range("C4").Value = 112.13
$ws.Range("D4").Value = 179.51
$ws.Range("E4").Value = 146.22
$ws.Range("F4").Value = 115.06
$ws.Range("G4").Value = 67.70999999999999
$ws.Range("H4").Value = 101.56
$ws.Range("I4").Value = 135.85
$ws.Range("J4").Value = 58.11
$ws.Range("K4").Value = 101.56
$ws.Range("L4").Value = 123.61

# Row 5: all columns D5:L5 = 6440 (same as existing C5)
$ws.Range("D5").Value = 6440
$ws.Range("E5").Value = 6440
$ws.Range("F5").Value = 6440
$ws.Range("G5").Value = 6440
$ws.Range("H5").Value = 6440
$ws.Range("I5").Value = 6440
$ws.Range("J5").Value = 6440
$ws.Range("K5").Value = 6440
$ws.Range("L5").Value = 6440

# Row 6: D6:L6 = 50 (C6 already 50, unchanged)
$ws.Range("D6").Value = 50
$ws.Range("E6").Value = 50
$ws.Range("F6").Value = 50
$ws.Range("G6").Value = 50
$ws.Range("H6").Value = 50
$ws.Range("I6").Value = 50
$ws.Range("J6").Value = 50
$ws.Range("K6").Value = 50
$ws.Range("L6").Value = 50
